# Weekly update: add a new week of "Zanahoria" (Mercado Mayorista Lo Valledor
# de Santiago) price observations at the top of the data block that starts at
# row 1272, pushing the existing rows down by 5 (dimension grows from
# A1:R1294 to A1:R1299).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before row 1272; everything that used to live in
# 1272:1294 shifts down to 1277:1299 automatically.
$ws.Rows("1272:1276").Insert()

# Seed the 5 new rows with the same "constant" column values (market id,
# market name, region, codreg, category id/name, variety, unit, kg/unit,
# classification) as the surrounding rows, using the row now sitting at
# 1277 (the old row 1272) as the template.
$ws.Range("A1277:R1277").Copy()
for ($r = 1272; $r -le 1276; $r++) {
    $ws.Range("A$($r):R$($r)").PasteSpecial()
}

# --- Row 1272: Camote ---
$ws.Range("D1272").Value = 44628
$ws.Range("I1272").Value = "Camote"
$ws.Range("J1272").Value = 450
$ws.Range("K1272").Value = 8000
$ws.Range("L1272").Value = 8000
$ws.Range("M1272").Value = 8000
$ws.Range("O1272").Value = "Chillán"
$ws.Range("P1272").Value = 400

# --- Row 1273: Primera ---
$ws.Range("D1273").Value = 44628
$ws.Range("I1273").Value = "Primera"
$ws.Range("J1273").Value = 1370
$ws.Range("K1273").Value = 9000
$ws.Range("L1273").Value = 9000
$ws.Range("M1273").Value = 9000
$ws.Range("O1273").Value = "Chillán"
$ws.Range("P1273").Value = 450

# --- Row 1274: Primera ---
$ws.Range("D1274").Value = 44628
$ws.Range("I1274").Value = "Primera"
$ws.Range("J1274").Value = 900
$ws.Range("K1274").Value = 8000
$ws.Range("L1274").Value = 8000
$ws.Range("M1274").Value = 8000
$ws.Range("O1274").Value = "Región Metropolitana"
$ws.Range("P1274").Value = 400

# --- Row 1275: Segunda ---
$ws.Range("D1275").Value = 44628
$ws.Range("I1275").Value = "Segunda"
$ws.Range("J1275").Value = 580
$ws.Range("K1275").Value = 7000
$ws.Range("L1275").Value = 7000
$ws.Range("M1275").Value = 7000
$ws.Range("O1275").Value = "Chillán"
$ws.Range("P1275").Value = 350

# --- Row 1276: Segunda ---
$ws.Range("D1276").Value = 44628
$ws.Range("I1276").Value = "Segunda"
$ws.Range("J1276").Value = 600
$ws.Range("K1276").Value = 7000
$ws.Range("L1276").Value = 7000
$ws.Range("M1276").Value = 7000
$ws.Range("O1276").Value = "Región Metropolitana"
$ws.Range("P1276").Value = 350

Write-Output "Inserted 5 new rows (1272:1276); sheet now spans $($ws.UsedRange.Address)"
